$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.797.46"
$ws.Range("E2").Value = "  +2.12%  "
$ws.Range("D3").Value = "1.858.09"
$ws.Range("E3").Value = "  +1.68%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'245.21"
$ws.Range("E5").Value = "  +1.38%  "
$ws.Range("D6").Value = "'0.6409"
$ws.Range("E6").Value = "  +3.20%  "
$ws.Range("D7").Value = "'0.9999"
$ws.Range("E7").Value = "  -0.15%  "
$ws.Range("D8").Value = "'0.07499"
$ws.Range("E8").Value = "  +2.22%  "
$ws.Range("D9").Value = "'0.2982"
$ws.Range("E9").Value = "  +2.88%  "
$ws.Range("D10").Value = "'24.15"
$ws.Range("E10").Value = "  +5.04%  "
$ws.Range("D11").Value = "'0.07682"
$ws.Range("E11").Value = "  -0.09%  "
$ws.Range("D12").Value = "1.877.03"
$ws.Range("E12").Value = "  +2.96%  "
$ws.Range("D13").Value = "'5.063"
$ws.Range("E13").Value = "  +2.02%  "
$ws.Range("D14").Value = "'0.6851"
$ws.Range("E14").Value = "  +3.02%  "
$ws.Range("D15").Value = "'84.05"
$ws.Range("E15").Value = "  +1.95%  "
$ws.Range("D16").Value = "'0.000009471"
$ws.Range("E16").Value = "  +5.80%  "
$ws.Range("D17").Value = "'6.061"
$ws.Range("E17").Value = "  +3.52%  "
$ws.Range("D18").Value = "29.762.34"
$ws.Range("E18").Value = "  +2.08%  "
$ws.Range("D19").Value = "2.116.32"
$ws.Range("E19").Value = "  +2.34%  "
$ws.Range("D20").Value = "'239.64"
$ws.Range("E20").Value = "  +0.36%  "
$ws.Range("D21").Value = "'12.69"
$ws.Range("E21").Value = "  +2.10%  "
$ws.Range("D22").Value = "'0.9998"
$ws.Range("E22").Value = "  -0.12%  "
$ws.Range("D23").Value = "'7.428"
$ws.Range("E23").Value = "  +1.62%  "
$ws.Range("E24").Value = "  -0.06%  "
$ws.Range("D25").Value = "'158.82"
$ws.Range("D26").Value = "'0.1430"
$ws.Range("E26").Value = "  +0.48%  "
$ws.Range("D27").Value = "'8.531"
$ws.Range("E27").Value = "  +0.65%  "
$ws.Range("D28").Value = "'17.95"
$ws.Range("E28").Value = "  +1.56%  "
$ws.Range("D29").Value = "'0.06223"
$ws.Range("E29").Value = "  +11.45%  "
$ws.Range("D30").Value = "'1.500"
$ws.Range("E30").Value = "  +1.11%  "
$ws.Range("D31").Value = "'1.273"
$ws.Range("E31").Value = "  +5.48%  "
$ws.Range("D32").Value = "'4.148"
$ws.Range("E32").Value = "  +1.24%  "
$ws.Range("D33").Value = "'4.118"
$ws.Range("E33").Value = "  +0.66%  "
$ws.Range("D34").Value = "'1.882"
$ws.Range("E34").Value = "  +2.08%  "
$ws.Range("D35").Value = "'1.161"
$ws.Range("E35").Value = "  +2.40%  "
$ws.Range("D36").Value = "'0.7324"
$ws.Range("E36").Value = "  -0.33%  "
$ws.Range("E37").Value = "  -0.95%  "
$ws.Range("D38").Value = "'2.857"
$ws.Range("E38").Value = "  +0.46%  "
$ws.Range("E39").Value = "  +1.70%  "
$ws.Range("D40").Value = "1.213.79"
$ws.Range("E40").Value = "  +0.22%  "
$ws.Range("D41").Value = "'0.9257"
$ws.Range("E41").Value = "  +0.84%  "
$ws.Range("D42").Value = "'6.173"
$ws.Range("E42").Value = "  -2.11%  "
$ws.Range("B43").Value = "PaxDollar"
$ws.Range("C43").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D43").Value = "'1.000"
$ws.Range("E43").Value = "  -0.15%  "
$ws.Range("B44").Value = "RocketPoolETH"
$ws.Range("C44").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D44").Value = "2.019.65"
$ws.Range("E44").Value = "  +2.48%  "
$ws.Range("D45").Value = "'102.04"
$ws.Range("E45").Value = "  +0.24%  "
$ws.Range("D46").Value = "'66.52"
$ws.Range("E46").Value = "  +2.56%  "
$ws.Range("E47").Value = "  +1.99%  "
$ws.Range("D48").Value = "'9.315"
$ws.Range("E48").Value = "  +1.65%  "
$ws.Range("D49").Value = "'0.4085"
$ws.Range("E49").Value = "  +1.43%  "
$ws.Range("D50").Value = "'0.1133"
$ws.Range("E50").Value = "  +1.71%  "
$ws.Range("D51").Value = "'0.05799"
$ws.Range("E51").Value = "  +0.68%  "